$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")

# Row 45: add C45 = "test"
$ws.Range("C45").Value = "test"

# Row 47 (new): B47 = "test tool", C47 = "tts", D47 = "description"
$ws.Range("B47").Value = "test tool"
$ws.Range("C47").Value = "tts"
$ws.Range("D47").Value = "description"

# Row 49 (new): B49 = "tool test", C49 = "ttt", D49 = "description"
$ws.Range("B49").Value = "tool test"
$ws.Range("C49").Value = "ttt"
$ws.Range("D49").Value = "description"

# Row 51 (new): B51 = "duplicate test tool", C51 = "dtt", D51 = "Presenting Creative Solutions Tool Kit"
$ws.Range("B51").Value = "duplicate test tool"
$ws.Range("C51").Value = "dtt"
$ws.Range("D51").Value = "Presenting Creative Solutions Tool Kit"

# Update the sheet view to match the saved selection/scroll position
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D51").Select()
